$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (single dot, valid float).
# These must be forced to remain text, matching the original inlineStr/shared-string
# cell type in the workbook (t="inlineStr" in source OOXML) instead of being
# auto-coerced to a numeric cell by Excel.
$textForceCells = @{
    'D5' = '248.76'
    'D6' = '0.667'
    'D7' = '58.89'
    'D12' = '15.88'
    'D14' = '0.839'
    'D15' = '5.79'
    'D17' = '18.17'
    'D19' = '75.29'
    'D21' = '5.44'
    'D22' = '238.53'
    'D26' = '169.48'
    'D27' = '9.44'
    'D28' = '20.19'
    'D30' = '4.86'
    'D31' = '1.12'
    'D32' = '0.0624'
    'D34' = '0.0911'
    'D38' = '1.35'
    'D39' = '0.107'
    'D40' = '3.15'
    'D41' = '5.10'
    'D43' = '17.47'
    'D45' = '96.81'
    'D46' = '2.49'
    'D50' = '3.68'
}

foreach ($ref in $textForceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textForceCells[$ref]
    $cell.ClearFormats()
}

# Cells whose new values are plain text (coin names, URLs, percentage strings,
# or price strings containing more than one "." which Excel cannot parse as a
# number anyway) can be assigned directly.
$plainCells = @{
    'D2' = '37.167.31'
    'E2' = '  -0.15%  '
    'D3' = '2.054.78'
    'E3' = '  -0.97%  '
    'E4' = '  -0.28%  '
    'E6' = '  -1.06%  '
    'E7' = '  -5.64%  '
    'E8' = '  +0.02%  '
    'E9' = '  +0.32%  '
    'E10' = '  -1.47%  '
    'E11' = '  +0.79%  '
    'E12' = '  +0.92%  '
    'D13' = '2.355.44'
    'E13' = '  +0.44%  '
    'E14' = '  +1.50%  '
    'E15' = '  +6.95%  '
    'D16' = '2.055.92'
    'E16' = '  -1.19%  '
    'E17' = '  +20.16%  '
    'D18' = '37.212.35'
    'E18' = '  +0.11%  '
    'E19' = '  +0.60%  '
    'E20' = '  -2.46%  '
    'E21' = '  -0.54%  '
    'E22' = '  -0.64%  '
    'E23' = '  +0.01%  '
    'E24' = '  +2.63%  '
    'E25' = '  +6.04%  '
    'E26' = '  -1.18%  '
    'E27' = '  +1.88%  '
    'E28' = '  -1.01%  '
    'E29' = '  -0.03%  '
    'E30' = '  +2.64%  '
    'E31' = '  +2.78%  '
    'E32' = '  -2.09%  '
    'E33' = '  +2.97%  '
    'E34' = '  +2.33%  '
    'E35' = '  -0.18%  '
    'E36' = '  -0.66%  '
    'E37' = '  -0.36%  '
    'B38' = 'TrustWalletToken'
    'C38' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E38' = '  -0.50%  '
    'B39' = 'Cronos'
    'C39' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E39' = '  -1.56%  '
    'E40' = '  +12.46%  '
    'E41' = '  +12.00%  '
    'E42' = '  -1.69%  '
    'E43' = '  -7.32%  '
    'E44' = '  -1.11%  '
    'E45' = '  -1.77%  '
    'E46' = '  -3.20%  '
    'D47' = '1.289.34'
    'E47' = '  -1.47%  '
    'E48' = '  -0.69%  '
    'E49' = '  -0.92%  '
    'E50' = '  -16.05%  '
    'D51' = '2.248.11'
    'E51' = '  -0.36%  '
}

foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
